$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$np1 = $s1.NotesPage
$np1.Shapes.Item(2).TextFrame.TextRange.Text = '各位同學，早晨！喺我哋上一次嘅講座入面，我哋討論咗牛頓運動定律。今日，我哋會透過介紹功、能量同功率嚟擴展呢啲概念。呢啲概念對於理解力點樣喺距離同時間上影響物體係至關重要嘅。'

$s2 = $p.Slides.Item(2)
$np2 = $s2.NotesPage
$np2.Shapes.Item(2).TextFrame.TextRange.Text = '功，在物理學中，被定義為一個力導致位移。如果一個力導致物體產生位移，那麼這個力就對物體做功。功的量被計算為力乘以距離，特別是力在位移方向上的分量。'

$s3 = $p.Slides.Item(3)
$np3 = $s3.NotesPage
$np3.Shapes.Item(2).TextFrame.TextRange.Text = '能量係做功嘅能力。能量有好多種形式，例如動能（郁動嘅能量）同勢能（儲存嘅能量）。能量守恆定律指出，能量唔能夠被創造或者被毀滅，只係由一種形式轉化成另一種形式。'

$s4 = $p.Slides.Item(4)
$np4 = $s4.NotesPage
$np4.Shapes.Item(2).TextFrame.TextRange.Text = '動能係物體因為郁動而擁有嘅能量。佢同物體嘅質量同埋佢速度嘅平方成正比。所以，一個郁得越快或者越重嘅物體，就會有越多嘅動能。'

$s5 = $p.Slides.Item(5)
$np5 = $s5.NotesPage
$np5.Shapes.Item(2).TextFrame.TextRange.Text = '勢能係儲存嘅能量。重力勢能取決於物體嘅質量、高度同地心吸力加速度。彈性勢能儲存喺拉伸或壓縮嘅彈性物料入面。'

$s6 = $p.Slides.Item(6)
$np6 = $s6.NotesPage
$np6.Shapes.Item(2).TextFrame.TextRange.Text = '最後，功率就係做功或者能量轉移嘅速率。佢話俾我哋知做功嘅速度有幾快。功率可以計做功除以時間，或者力乘以速度。我哋關於功、能量同功率嘅講座到此結束。請為下次課堂嘅一個小測驗做準備。'
